$d = $word.ActiveDocument

function Set-ParaText($paraIndex, $newText) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    # Exclude the trailing paragraph mark from the range, wipe any existing
    # runs/line-breaks, then insert the replacement as one fresh run. This
    # avoids Find/Replace's smart-quote autocorrect and avoids Range.Text's
    # habit of keeping a stale xml:space="preserve" from the old run.
    $r.MoveEnd(1, -1) | Out-Null
    $r.Text = ""
    $r.InsertBefore($newText)
}

# --- Paragraph 1: date + title (two separate <w:t> runs joined by a <w:br/>) ---
$p1 = $d.Paragraphs(1).Range
$p1.Find.Execute("-14.10.24", $true, $false, $false, $false, $false, $true, 1, $false, "-12.10.24", 2) | Out-Null
$p1.Find.Execute("LLMS KNOW MORE THAN THEY SHOW: ON THE IN-TRINSIC REPRESENTATION OF LLM HALLUCINATIONS", $true, $false, $false, $false, $false, $true, 1, $false, "GSM-Symbolic: Understanding the Limitations of Mathematical Reasoning in Large Language Models", 2) | Out-Null

# --- Paragraph 2: summary paragraph, full replacement ---
Set-ParaText 2 "האם מודלי שפה גדולים מסוגלים לעשות ריזונינג? השאלה הזו מעסיקה חוקרים רבים לאור יכולות די מרשימות שמודלי שפה מפגינים בפתרון שאלות לא פשוטות (אבל רק בתנאים מסוימים 🙂). המאמר בוחן את יכולות החשיבה המתמטית(שזה תת-יכולת של ריזונינג כללי) של LLMs ומציג את GSM-Symbolic, בנצ'מרק חדש לבחינת יכולות אלו שהם פיתחו."

# --- Paragraph 3: was " המחברים חקרו..." + <br><br> + "שיפור זיהוי שגיאות:" -> single new paragraph ---
Set-ParaText 3 "החוקרים מצאו שביצועי LLMs(נבחן מגוון רחב של מודלים: Gemma, Phi, Mistral, Llama3, GPT-4 ו-o1) משתנים באופן משמעותי כאשר משנים מעט את השאלות המתמטיות, מה שמעלה ספקות לגבי אמינות המדדים הקיימים. הביצועים של רוב המודלים יורדים כאשר עוברים מ-GSM8K המקורי ל-GSM-Symbolic, מה שמרמז על אפשרות של זיהום דאטה (contamination) במהלך האימון (כתבתי על זה לא מעט)."

# --- Paragraph 4: was "המחברים גילו..." + <br><br> + "הכללה בין משימות:" -> single new paragraph ---
Set-ParaText 4 "בנוסף המודלים מראים רגישות גבוהה יותר לשינויים במספרים מאשר לשינויים בשמות עצם, מה שמעיד על חוסר יציבות ביכולות החשיבה שלהם. ככל שמספר המשפטים בשאלה עולה, הביצועים יורדים והשונות בביצועים עולה, מה שמצביע על קושי בטיפול בשאלות מורכבות יותר."

# --- Paragraph 5 ---
Set-ParaText 5 "החוקרים יצרו בנצ'מארק GSM-NoOp, שבו נוספו משפטים לא רלוונטיים לשאלות, וגילו ירידה דרמטית בביצועים של כל המודלים. אפילו כאשר ניתנו למודלים דוגמאות של אותה שאלה או שאלות דומות, הם התקשו להתגבר על האתגרים של GSM-NoOp. "

# --- Paragraph 6 ---
Set-ParaText 6 "המחקר מצא שאימון נוסף על משימות קלות יותר וגם הגדלת כמות דאטה לאימון לא שיפרו את הביצועים במשימות מורכבות יותר."

# --- Paragraph 7 ---
Set-ParaText 7 "קצת מנחם שלפחות מודלים חדשים יותר, כמו o1-preview ו-o1-mini, הראו ביצועים חזקים יותר, אך עדיין סבלו מהמגבלות שזוהו במחקר"

# --- Paragraph 8 ---
Set-ParaText 8 "הממצאים מעלים ספקות לגבי היכולת האמיתית של LLMs לבצע חשיבה מתמטית פורמלית. נראה כי המודלים מסתמכים יותר על התאמת תבניות מאשר על חשיבה לוגית אמיתית. המחקר מדגיש את הצורך בשיטות הערכה אמינות יותר ובמחקר נוסף על יכולות החשיבה של מודלי שפה גדולים."

# --- Paragraph 9: replaced in-place with the new arXiv URL ---
Set-ParaText 9 "https://arxiv.org/abs/2410.05229"

# --- Paragraphs 10 and 11: removed entirely (old "ממצא זה מרמז..." paragraph and the old URL paragraph) ---
$d.Paragraphs(11).Range.Delete() | Out-Null
$d.Paragraphs(10).Range.Delete() | Out-Null
